$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 - who/when PyTorch was developed ---
$r1 = $d.Content
$null = $r1.Find.Execute(
    " was primarily developed by Facebook's AI Research lab (FAIR), which is now part of Meta AI.^sIts development began in 2016.^s",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " was originally developed by Facebook AI Research (FAIR), which is now known as Meta AI.", 2)

# --- Edit 2: paragraph 2, opening clause ---
$r2 = $d.Content
$null = $r2.Find.Execute(
    "While Meta AI was the primary developer and maintainer for a significant period, the administration of the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The development of ", 2)

# --- Edit 3: paragraph 2, middle clause ---
$r3 = $d.Content
$null = $r3.Find.Execute(
    " project was later handed over to the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " began in 2016 as an open-source machine learning framework, building upon the Torch library. While initially a project within Facebook, its administration was transferred to the neutral ", 2)

# --- Edit 4: paragraph 2, closing clause (spans the final "PyTorch" + "." runs) ---
$r4 = $d.Content
$null = $r4.Find.Execute(
    " Foundation in 2022. This foundation operates under the umbrella of the Linux Foundation and is responsible for coordinating the future development and ecosystem of PyTorch.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Foundation, part of the Linux Foundation, in 2022 to foster broader community collaboration and development.", 2)

# --- Edit 5: remove the now-dropped "original authors" paragraph entirely ---
$d.Paragraphs.Item(3).Range.Delete()
